$d = $word.ActiveDocument

# Title (Heading1) and the matching bold run near the end share the same old/new text.
$d.Paragraphs.Item(1).Range.Find.Execute("Play Aztec Wilds Free: A Thrilling Game with Multiple Wilds", $true, $false, $false, $false, $false, $true, 1, $false, "Play Aztec Wilds Free - Review of Gameplay, Wilds, Symbols, and Volatility", 2) | Out-Null

# "What we like" bullet list.
$d.Paragraphs.Item(42).Range.Find.Execute("1,024 possible ways to win", $true, $false, $false, $false, $false, $true, 1, $false, "Highly polished graphics", 2) | Out-Null
$d.Paragraphs.Item(44).Range.Find.Execute("Highly polished graphics", $true, $false, $false, $false, $false, $true, 1, $false, "Multiple Wild symbols", 2) | Out-Null
$d.Paragraphs.Item(45).Range.Find.Execute("A variety of Wild Symbols", $true, $false, $false, $false, $false, $true, 1, $false, "Cascading symbols and Free Spins", 2) | Out-Null

# "What we don't like" bullet list. Use a plain Find (no auto-replace) followed by a direct
# Range.Text assignment so the straight apostrophes in "'sticky'" are not mangled into curly
# smart quotes by the runtime's Find-and-replace auto-formatting.
$r47 = $d.Paragraphs.Item(47).Range
$r47.Find.Execute("High volatility") | Out-Null
$r47.Text = "Limited appearance of 'sticky' Wild"

$d.Paragraphs.Item(48).Range.Find.Execute("No bonus games", $true, $false, $false, $false, $false, $true, 1, $false, "High volatility", 2) | Out-Null

# Bold title recap near the end of the document.
$d.Paragraphs.Item(49).Range.Find.Execute("Play Aztec Wilds Free: A Thrilling Game with Multiple Wilds", $true, $false, $false, $false, $false, $true, 1, $false, "Play Aztec Wilds Free - Review of Gameplay, Wilds, Symbols, and Volatility", 2) | Out-Null

# Italic meta description at the very end.
$d.Paragraphs.Item(50).Range.Find.Execute("Discover a new world of slot games with Aztec Wilds free play version. A game with a unique gameplay system filled with wild symbols and Free Spins.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Aztec Wilds gameplay, Wild symbols, Cascading symbols, and volatility. Play for free!", 2) | Out-Null
